$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column F (etaSqp) width slightly for display
$ws.Range("F1").EntireColumn.ColumnWidth = 12.7109375

# Updated ANOVA statistics
$ws.Range("D2").Value = 375858.88519869605
$ws.Range("F2").Value = 0.99869801660302171

$ws.Range("D3").Value = 4.009086477938423
$ws.Range("E3").Value = 0.0458062517406711
$ws.Range("F3").Value = 0.0081154103996778654

$ws.Range("D4").Value = 945.88974503464829
$ws.Range("F4").Value = 0.65874817220860071

$ws.Range("D5").Value = 138.89344621813584
$ws.Range("F5").Value = 0.22085370272718619
